# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45181 (2023-09-12) to 45182 (2023-09-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 360; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
